# Add "Flaxon 650" / "flaxon 650" drug entries to the Drugs list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Flaxon 650" just before "Flucloxacillin" (currently row 30) ---
$ws.Rows.Item(30).Insert()
$ws.Cells.Item(30, 2).Value = "Flaxon 650"

# --- Insert an extra "Paracetamol500" row just before "Penicillin G (Benzylpenicillin)" ---
# (after the row above was inserted, "Penicillin G ..." is now at row 52)
$ws.Rows.Item(52).Insert()
$ws.Cells.Item(52, 2).Value = "Paracetamol500"

# --- Append "flaxon 650" as a new last row ---
$ws.Cells.Item(70, 2).Value = "flaxon 650"

# --- Re-number the Index column (A2:A70) sequentially: 1, 2, 3, ... ---
for ($r = 2; $r -le 70; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
